$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-20 Tuesday", $true, $true, $false, $false, $false, $true, 1, $false, "2026-01-21 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("225×9=", $true, $true, $false, $false, $false, $true, 1, $false, "294×6=", 2) | Out-Null
$d.Content.Find.Execute("812×9=", $true, $true, $false, $false, $false, $true, 1, $false, "564×2=", 2) | Out-Null
$d.Content.Find.Execute("950×7=", $true, $true, $false, $false, $false, $true, 1, $false, "796×4=", 2) | Out-Null
$d.Content.Find.Execute("339×7=", $true, $true, $false, $false, $false, $true, 1, $false, "238×7=", 2) | Out-Null
$d.Content.Find.Execute("623×8=", $true, $true, $false, $false, $false, $true, 1, $false, "698×7=", 2) | Out-Null
$d.Content.Find.Execute("173×6=", $true, $true, $false, $false, $false, $true, 1, $false, "208×4=", 2) | Out-Null
$d.Content.Find.Execute("528×9=", $true, $true, $false, $false, $false, $true, 1, $false, "802×7=", 2) | Out-Null
$d.Content.Find.Execute("791×6=", $true, $true, $false, $false, $false, $true, 1, $false, "690×2=", 2) | Out-Null
$d.Content.Find.Execute("641×4=", $true, $true, $false, $false, $false, $true, 1, $false, "714×9=", 2) | Out-Null
$d.Content.Find.Execute("579×4=", $true, $true, $false, $false, $false, $true, 1, $false, "135×7=", 2) | Out-Null
$d.Content.Find.Execute("676×4=", $true, $true, $false, $false, $false, $true, 1, $false, "306×5=", 2) | Out-Null
$d.Content.Find.Execute("379×4=", $true, $true, $false, $false, $false, $true, 1, $false, "328×9=", 2) | Out-Null
$d.Content.Find.Execute("209×5=", $true, $true, $false, $false, $false, $true, 1, $false, "481×6=", 2) | Out-Null
$d.Content.Find.Execute("437×9=", $true, $true, $false, $false, $false, $true, 1, $false, "404×2=", 2) | Out-Null
$d.Content.Find.Execute("989×7=", $true, $true, $false, $false, $false, $true, 1, $false, "251×4=", 2) | Out-Null
$d.Content.Find.Execute("929×8=", $true, $true, $false, $false, $false, $true, 1, $false, "757×4=", 2) | Out-Null
$d.Content.Find.Execute("191×3=", $true, $true, $false, $false, $false, $true, 1, $false, "948×8=", 2) | Out-Null
$d.Content.Find.Execute("666×6=", $true, $true, $false, $false, $false, $true, 1, $false, "372×6=", 2) | Out-Null
$d.Content.Find.Execute("451×3=", $true, $true, $false, $false, $false, $true, 1, $false, "730×5=", 2) | Out-Null
$d.Content.Find.Execute("995×3=", $true, $true, $false, $false, $false, $true, 1, $false, "525×2=", 2) | Out-Null
$d.Content.Find.Execute("783×2=", $true, $true, $false, $false, $false, $true, 1, $false, "914×7=", 2) | Out-Null
$d.Content.Find.Execute("713×4=", $true, $true, $false, $false, $false, $true, 1, $false, "334×7=", 2) | Out-Null
$d.Content.Find.Execute("957×9=", $true, $true, $false, $false, $false, $true, 1, $false, "700×2=", 2) | Out-Null
$d.Content.Find.Execute("345×8=", $true, $true, $false, $false, $false, $true, 1, $false, "873×9=", 2) | Out-Null
$d.Content.Find.Execute("636×6=", $true, $true, $false, $false, $false, $true, 1, $false, "264×9=", 2) | Out-Null

Write-Output "Done"